# Auto-generated edit script applying scheduled price-data refresh
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 2156.4614  # H4
$ws.Cells.Item(4, 9).Value = 1550.3334  # I4
$ws.Cells.Item(4, 10).Value = 3520.25  # J4
$ws.Cells.Item(4, 11).Value = 1550.3334  # K4
$ws.Cells.Item(4, 12).Value = 3520.25  # L4
$ws.Cells.Item(4, 13).Value = -1436.3334  # M4
$ws.Cells.Item(4, 14).Value = -3748.25  # N4
$ws.Cells.Item(42, 8).Value = 2525.625  # H42
$ws.Cells.Item(42, 9).Value = 381.9  # I42
$ws.Cells.Item(42, 10).Value = 6098.5  # J42
$ws.Cells.Item(42, 11).Value = 1145.7  # K42
$ws.Cells.Item(42, 12).Value = 18295.5  # L42
$ws.Cells.Item(42, 13).Value = -915.6999999999998  # M42
$ws.Cells.Item(42, 14).Value = -18755.5  # N42
$ws.Cells.Item(62, 8).Value = 7944.375  # H62
$ws.Cells.Item(62, 9).Value = 7897.6665  # I62
$ws.Cells.Item(62, 10).Value = 7972.4  # J62
$ws.Cells.Item(62, 11).Value = 7897.6665  # K62
$ws.Cells.Item(62, 12).Value = 7972.4  # L62
$ws.Cells.Item(62, 13).Value = -7273.6665  # M62
$ws.Cells.Item(62, 14).Value = -9220.4  # N62
$ws.Cells.Item(65, 8).Value = 7944.375  # H65
$ws.Cells.Item(65, 9).Value = 7897.6665  # I65
$ws.Cells.Item(65, 10).Value = 7972.4  # J65
$ws.Cells.Item(65, 11).Value = 39488.3325  # K65
$ws.Cells.Item(65, 12).Value = 39862  # L65
$ws.Cells.Item(65, 13).Value = -36368.3325  # M65
$ws.Cells.Item(65, 14).Value = -46102  # N65
$ws.Cells.Item(112, 8).Value = 2899  # H112

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 289.81818  # H2
$ws.Cells.Item(2, 9).Value = 289.81818  # I2
$ws.Cells.Item(2, 11).Value = 289.81818  # K2
$ws.Cells.Item(2, 13).Value = -176.81818  # M2
$ws.Cells.Item(5, 8).Value = 717.25  # H5
$ws.Cells.Item(5, 9).Value = 731.3333  # I5
$ws.Cells.Item(5, 11).Value = 731.3333  # K5
$ws.Cells.Item(5, 13).Value = -619.3333  # M5
$ws.Cells.Item(45, 8).Value = 4900  # H45
$ws.Cells.Item(45, 9).Value = 4900  # I45
$ws.Cells.Item(45, 11).Value = 4900  # K45
$ws.Cells.Item(45, 13).Value = -4523  # M45
$ws.Cells.Item(46, 8).Value = 7269.3335  # H46
$ws.Cells.Item(46, 10).Value = 7269.3335  # J46
$ws.Cells.Item(46, 12).Value = 7269.3335  # L46
$ws.Cells.Item(46, 14).Value = -7907.3335  # N46
$ws.Cells.Item(61, 8).Value = 4509.1  # H61
$ws.Cells.Item(61, 9).Value = 1958.8572  # I61
$ws.Cells.Item(61, 10).Value = 5882.3076  # J61
$ws.Cells.Item(61, 11).Value = 1958.8572  # K61
$ws.Cells.Item(61, 12).Value = 5882.3076  # L61
$ws.Cells.Item(61, 13).Value = -1746.8572  # M61
$ws.Cells.Item(61, 14).Value = -6306.3076  # N61
$ws.Cells.Item(74, 8).Value = 2593.5  # H74
$ws.Cells.Item(74, 9).Value = 1967.5385  # I74
$ws.Cells.Item(74, 11).Value = 1967.5385  # K74
$ws.Cells.Item(74, 13).Value = -1093.5385  # M74
$ws.Cells.Item(77, 8).Value = 2593.5  # H77
$ws.Cells.Item(77, 9).Value = 1967.5385  # I77
$ws.Cells.Item(77, 11).Value = 9837.692500000001  # K77
$ws.Cells.Item(77, 13).Value = -5469.692500000001  # M77
$ws.Cells.Item(102, 8).Value = 1292  # H102
$ws.Cells.Item(102, 9).Value = 1292  # I102
$ws.Cells.Item(102, 11).Value = 1292  # K102
$ws.Cells.Item(102, 13).Value = 330  # M102
$ws.Cells.Item(113, 8).Value = 125999.57  # H113
$ws.Cells.Item(113, 10).Value = 125999.57  # J113
$ws.Cells.Item(113, 12).Value = 125999.57  # L113
$ws.Cells.Item(113, 14).Value = -134677.57  # N113
$ws.Cells.Item(116, 8).Value = 289.81818  # H116
$ws.Cells.Item(116, 9).Value = 289.81818  # I116
$ws.Cells.Item(116, 11).Value = 289.81818  # K116
$ws.Cells.Item(116, 13).Value = 2004.18182  # M116
$ws.Cells.Item(122, 8).Value = 1142.2  # H122
$ws.Cells.Item(122, 9).Value = 1149.25  # I122
$ws.Cells.Item(122, 10).Value = 1114  # J122
$ws.Cells.Item(122, 11).Value = 3447.75  # K122
$ws.Cells.Item(122, 12).Value = 3342  # L122
$ws.Cells.Item(122, 13).Value = -997.75  # M122
$ws.Cells.Item(122, 14).Value = -8242  # N122
$ws.Cells.Item(132, 8).Value = 1445.6875  # H132
$ws.Cells.Item(132, 9).Value = 1445.6875  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 4337.0625  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).Value = -1807.0625  # M132
$ws.Cells.Item(132, 14).ClearContents()  # N132 (removed)
$ws.Cells.Item(136, 8).Value = 4509.1  # H136
$ws.Cells.Item(136, 9).Value = 1958.8572  # I136
$ws.Cells.Item(136, 10).Value = 5882.3076  # J136
$ws.Cells.Item(136, 11).Value = 5876.571599999999  # K136
$ws.Cells.Item(136, 12).Value = 17646.9228  # L136
$ws.Cells.Item(136, 13).Value = -3326.571599999999  # M136
$ws.Cells.Item(136, 14).Value = -22746.9228  # N136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 289.81818  # H3
$ws.Cells.Item(3, 9).Value = 289.81818  # I3
$ws.Cells.Item(3, 11).Value = 289.81818  # K3
$ws.Cells.Item(3, 13).Value = -175.81818  # M3
$ws.Cells.Item(4, 8).Value = 717.25  # H4
$ws.Cells.Item(4, 9).Value = 731.3333  # I4
$ws.Cells.Item(4, 11).Value = 731.3333  # K4
$ws.Cells.Item(4, 13).Value = -616.3333  # M4
$ws.Cells.Item(22, 8).Value = 582.4  # H22
$ws.Cells.Item(22, 9).Value = 485.625  # I22
$ws.Cells.Item(22, 11).Value = 485.625  # K22
$ws.Cells.Item(22, 13).Value = -312.625  # M22
$ws.Cells.Item(36, 8).Value = 519  # H36
$ws.Cells.Item(36, 9).Value = 519  # I36
$ws.Cells.Item(36, 11).Value = 519  # K36
$ws.Cells.Item(36, 13).Value = 15  # M36

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1468.72  # H31
$ws.Cells.Item(31, 9).Value = 1246.6666  # I31
$ws.Cells.Item(31, 10).Value = 2634.5  # J31
$ws.Cells.Item(31, 11).Value = 1246.6666  # K31
$ws.Cells.Item(31, 12).Value = 2634.5  # L31
$ws.Cells.Item(31, 13).Value = -951.6666  # M31
$ws.Cells.Item(31, 14).Value = -3224.5  # N31
$ws.Cells.Item(34, 8).Value = 1468.72  # H34
$ws.Cells.Item(34, 9).Value = 1246.6666  # I34
$ws.Cells.Item(34, 10).Value = 2634.5  # J34
$ws.Cells.Item(34, 11).Value = 1246.6666  # K34
$ws.Cells.Item(34, 12).Value = 2634.5  # L34
$ws.Cells.Item(34, 13).Value = -1044.6666  # M34
$ws.Cells.Item(34, 14).Value = -3038.5  # N34
$ws.Cells.Item(99, 8).Value = 5625  # H99
$ws.Cells.Item(99, 9).Value = 1800  # I99
$ws.Cells.Item(99, 11).Value = 1800  # K99
$ws.Cells.Item(99, 13).Value = -302  # M99
$ws.Cells.Item(126, 8).Value = 5625  # H126
$ws.Cells.Item(126, 9).Value = 1800  # I126
$ws.Cells.Item(126, 11).Value = 5400  # K126
$ws.Cells.Item(126, 13).Value = -2930  # M126
$ws.Cells.Item(141, 8).Value = 750000  # H141
$ws.Cells.Item(141, 10).Value = 750000  # J141
$ws.Cells.Item(141, 12).Value = 750000  # L141
$ws.Cells.Item(141, 14).Value = -760360  # N141

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 37526788  # H11
$ws.Cells.Item(11, 9).Value = 41696344  # I11
$ws.Cells.Item(11, 11).Value = 125089032  # K11
$ws.Cells.Item(11, 13).Value = -125088892  # M11
$ws.Cells.Item(12, 8).Value = 122  # H12
$ws.Cells.Item(12, 9).Value = 5.25  # I12
$ws.Cells.Item(12, 10).Value = 215.4  # J12
$ws.Cells.Item(12, 11).Value = 15.75  # K12
$ws.Cells.Item(12, 12).Value = 646.2  # L12
$ws.Cells.Item(12, 13).Value = 157.25  # M12
$ws.Cells.Item(12, 14).Value = -992.2  # N12
$ws.Cells.Item(41, 8).Value = 2099  # H41
$ws.Cells.Item(41, 9).Value = 2099  # I41
$ws.Cells.Item(41, 10).Value = 0  # J41
$ws.Cells.Item(41, 11).Value = 6297  # K41
$ws.Cells.Item(41, 12).Value = 0  # L41
$ws.Cells.Item(41, 13).Value = -5959  # M41
$ws.Cells.Item(41, 14).ClearContents()  # N41 (removed)
$ws.Cells.Item(107, 8).Value = 458.4  # H107
$ws.Cells.Item(107, 10).Value = 696.6667  # J107
$ws.Cells.Item(107, 12).Value = 2090.0001  # L107
$ws.Cells.Item(107, 14).Value = -5930.0001  # N107
$ws.Cells.Item(137, 8).Value = 2280.3125  # H137
$ws.Cells.Item(137, 9).Value = 1730.875  # I137
$ws.Cells.Item(137, 10).Value = 2829.75  # J137
$ws.Cells.Item(137, 11).Value = 5192.625  # K137
$ws.Cells.Item(137, 12).Value = 8489.25  # L137
$ws.Cells.Item(137, 13).Value = -92.625  # M137
$ws.Cells.Item(137, 14).Value = -18689.25  # N137
$ws.Cells.Item(138, 8).Value = 3113.6  # H138
$ws.Cells.Item(138, 9).Value = 3113.6  # I138
$ws.Cells.Item(138, 11).Value = 9340.799999999999  # K138
$ws.Cells.Item(138, 13).Value = -4200.799999999999  # M138

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 344.5  # H9
$ws.Cells.Item(9, 9).Value = 344.5  # I9
$ws.Cells.Item(9, 11).Value = 344.5  # K9
$ws.Cells.Item(9, 13).Value = -174.5  # M9
$ws.Cells.Item(59, 8).Value = 0  # H59
$ws.Cells.Item(59, 10).Value = 0  # J59
$ws.Cells.Item(59, 12).Value = 0  # L59
$ws.Cells.Item(59, 14).ClearContents()  # N59 (removed)
$ws.Cells.Item(112, 8).Value = 150000  # H112
$ws.Cells.Item(112, 10).Value = 150000  # J112
$ws.Cells.Item(112, 12).Value = 150000  # L112
$ws.Cells.Item(112, 14).Value = -152216  # N112
$ws.Cells.Item(113, 8).Value = 1220.8572  # H113
$ws.Cells.Item(113, 9).Value = 1220.8572  # I113
$ws.Cells.Item(113, 11).Value = 1220.8572  # K113
$ws.Cells.Item(113, 13).Value = 949.1428000000001  # M113
$ws.Cells.Item(132, 8).Value = 783.63635  # H132
$ws.Cells.Item(132, 9).Value = 783.63635  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 2350.90905  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).Value = 179.0909499999998  # M132
$ws.Cells.Item(132, 14).ClearContents()  # N132 (removed)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(31, 8).Value = 868.5  # H31
$ws.Cells.Item(31, 9).Value = 810.6667  # I31
$ws.Cells.Item(31, 10).Value = 903.2  # J31
$ws.Cells.Item(31, 11).Value = 810.6667  # K31
$ws.Cells.Item(31, 12).Value = 903.2  # L31
$ws.Cells.Item(31, 13).Value = -562.6667  # M31
$ws.Cells.Item(31, 14).Value = -1399.2  # N31
$ws.Cells.Item(40, 8).Value = 4124.75  # H40
$ws.Cells.Item(40, 9).Value = 3499.6667  # I40
$ws.Cells.Item(40, 11).Value = 3499.6667  # K40
$ws.Cells.Item(40, 13).Value = -3363.6667  # M40
$ws.Cells.Item(55, 8).Value = 360.5  # H55
$ws.Cells.Item(55, 9).Value = 314.5  # I55
$ws.Cells.Item(55, 10).Value = 636.5  # J55
$ws.Cells.Item(55, 11).Value = 314.5  # K55
$ws.Cells.Item(55, 12).Value = 636.5  # L55
$ws.Cells.Item(55, 13).Value = -141.5  # M55
$ws.Cells.Item(55, 14).Value = -982.5  # N55
$ws.Cells.Item(122, 8).Value = 7353.643  # H122
$ws.Cells.Item(122, 9).Value = 7770.706  # I122
$ws.Cells.Item(122, 10).Value = 6709.091  # J122
$ws.Cells.Item(122, 11).Value = 23312.118  # K122
$ws.Cells.Item(122, 12).Value = 20127.273  # L122
$ws.Cells.Item(122, 13).Value = -20862.118  # M122
$ws.Cells.Item(122, 14).Value = -25027.273  # N122

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 29999  # H52
$ws.Cells.Item(52, 9).Value = 29999  # I52
$ws.Cells.Item(52, 10).Value = 0  # J52
$ws.Cells.Item(52, 11).Value = 29999  # K52
$ws.Cells.Item(52, 12).Value = 0  # L52
$ws.Cells.Item(52, 13).Value = -29773  # M52
$ws.Cells.Item(52, 14).ClearContents()  # N52 (removed)
$ws.Cells.Item(81, 8).Value = 1429569.9  # H81
$ws.Cells.Item(81, 9).Value = 1237.6  # I81
$ws.Cells.Item(81, 11).Value = 2475.2  # K81
$ws.Cells.Item(81, 13).Value = -1414.2  # M81
$ws.Cells.Item(84, 8).Value = 1429569.9  # H84
$ws.Cells.Item(84, 9).Value = 1237.6  # I84
$ws.Cells.Item(84, 11).Value = 12376  # K84
$ws.Cells.Item(84, 13).Value = -7072  # M84

